# Bug Fixed backstabbing logic and Modified Performance stats
#
# Updates the "Own + Selfish" performance numbers (Strat 1-4 rows) and the
# "All strats (Own + Selfish)" summary row. The three summary-row cells
# also lose their yellow run-level highlight (the paragraph-mark highlight
# on the cell stays untouched).

$d = $word.ActiveDocument

function Replace-Number($oldText, $newText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($oldText, $true, $false, $false, $false, $false, `
                  $true, 1, $false, $newText, 2) | Out-Null
}

function Replace-HighlightedNumber($oldText, $newText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Highlight = $false
    $find.Execute($oldText, $true, $false, $false, $false, $false, `
                  $true, 1, $false, $newText, 2, $true) | Out-Null
}

# Strat 1 (Own + Selfish)
Replace-Number "331.22222222222223" "421.6666666666667"
Replace-Number "222.11111111111111" "242.33333333333334"
Replace-Number "263.22222222222223" "314.8888888888889"

# Strat 2 (Own + Selfish)
Replace-Number "434.3333333333333" "494.6666666666667"
Replace-Number "454.77777777777777" "502.55555555555554"
Replace-Number "438.44444444444446" "449.6666666666667"

# Strat 3 (Own + Selfish)
Replace-Number "277.44444444444446" "360.22222222222223"
Replace-Number "257.3333333333333" "299.0"
Replace-Number "326.8888888888889" "400.3333333333333"

# Strat 4 (Own + Selfish)
Replace-Number "268.0" "316.44444444444446"
Replace-Number "288.3333333333333" "377.0"
Replace-Number "322.1111111111111" "354.55555555555554"

# All strats (Own + Selfish) - also drop the yellow run highlight
Replace-HighlightedNumber "404.8888888888889" "422.44444444444446"
Replace-HighlightedNumber "473.3333333333333" "440.22222222222223"
Replace-HighlightedNumber "440.6666666666667" "511.44444444444446"
